$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8-10 block
$ws.Range("A8").Value = "Amen Thompson"
$ws.Range("B8").Value = "SG,SF"
$ws.Range("C8").Value = "Houston Rockets"

$ws.Range("A9").Value = "Buddy Hield"
$ws.Range("B9").Value = "SG,SF"
$ws.Range("C9").Value = "Golden State Warriors"

$ws.Range("A10").Value = "Nikola Vucevic"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "Chicago Bulls"

# Rows 12-14 block
$ws.Range("A12").Value = "Mikal Bridges"
$ws.Range("B12").Value = "SG,SF,PF"
$ws.Range("C12").Value = "New York Knicks"

$ws.Range("A13").Value = "DeMar DeRozan"
$ws.Range("B13").Value = "SF,PF"
$ws.Range("C13").Value = "Sacramento Kings"

$ws.Range("A14").Value = "Jusuf Nurkic"
$ws.Range("B14").Value = "C"
$ws.Range("C14").Value = "Phoenix Suns"
